# Update the dSF column (column F) values for several rows as part of a
# "repull data, push all data, mean calculation" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 4
    12 = 1
    15 = 0
    26 = -5
    32 = -3
    36 = 7
    41 = -6
    47 = -4
    53 = -4
    62 = 1
    66 = 0
    69 = -4
    78 = 2
    86 = -5
    87 = 2
    93 = 1
    94 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
